$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force text type so numeric-looking strings (e.g. "240.04") are not
    # reinterpreted as numbers, matching the inline/shared-string cells in the source.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Rows 46 and 47 swap coin order (VeChain <-> Filecoin) with updated data
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D46") "5.64"
$ws.Range("E46").Value = "  +3.59%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0415"
$ws.Range("E47").Value = "  -1.78%  "

# Price (column D) and volume-change (column E) updates
Set-TextValue $ws.Range("D2") "95.208.60"
$ws.Range("E2").Value = "  -1.11%  "
Set-TextValue $ws.Range("D3") "3.451.06"
$ws.Range("E3").Value = "  +4.13%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "240.04"
$ws.Range("E5").Value = "  -3.42%  "
Set-TextValue $ws.Range("D6") "642.80"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  +6.24%  "
Set-TextValue $ws.Range("D8") "0.403"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("E9").Value = "  +0.00%  "
Set-TextValue $ws.Range("D10") "1.00"
$ws.Range("E10").Value = "  +2.07%  "
Set-TextValue $ws.Range("D11") "3.451.14"
$ws.Range("E11").Value = "  +4.21%  "
Set-TextValue $ws.Range("D12") "0.198"
$ws.Range("E12").Value = "  -3.58%  "
Set-TextValue $ws.Range("D13") "41.78"
$ws.Range("E13").Value = "  +4.82%  "
Set-TextValue $ws.Range("D14") "6.10"
$ws.Range("E14").Value = "  +1.28%  "
Set-TextValue $ws.Range("D15") "94.778.96"
$ws.Range("E15").Value = "  -1.30%  "
Set-TextValue $ws.Range("D16") "4.088.78"
$ws.Range("E16").Value = "  +3.75%  "
Set-TextValue $ws.Range("D17") "0.0000256"
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("E18").Value = "  +0.14%  "
Set-TextValue $ws.Range("D19") "3.412.66"
$ws.Range("E19").Value = "  +2.63%  "
Set-TextValue $ws.Range("D20") "17.78"
$ws.Range("E20").Value = "  +4.85%  "
Set-TextValue $ws.Range("D21") "11.41"
$ws.Range("E21").Value = "  +9.41%  "
$ws.Range("E22").Value = "  -4.23%  "
Set-TextValue $ws.Range("D23") "502.32"
$ws.Range("E23").Value = "  +0.52%  "
Set-TextValue $ws.Range("D24") "3.17"
$ws.Range("E24").Value = "  -4.74%  "
Set-TextValue $ws.Range("D25") "0.0000192"
$ws.Range("E25").Value = "  -1.99%  "
Set-TextValue $ws.Range("D26") "6.45"
$ws.Range("E26").Value = "  -1.16%  "
Set-TextValue $ws.Range("D27") "91.78"
$ws.Range("E27").Value = "  -4.01%  "
Set-TextValue $ws.Range("D28") "3.625.20"
$ws.Range("E28").Value = "  +3.44%  "
Set-TextValue $ws.Range("D29") "12.00"
$ws.Range("E29").Value = "  +0.55%  "
Set-TextValue $ws.Range("D30") "11.73"
$ws.Range("E30").Value = "  +7.20%  "
$ws.Range("E31").Value = "  +0.12%  "
Set-TextValue $ws.Range("D32") "2.75"
$ws.Range("E32").Value = "  +11.85%  "
Set-TextValue $ws.Range("D33") "0.137"
$ws.Range("E33").Value = "  -3.39%  "
Set-TextValue $ws.Range("D34") "0.184"
$ws.Range("E34").Value = "  -0.42%  "
Set-TextValue $ws.Range("D35") "30.97"
$ws.Range("E35").Value = "  +11.30%  "
Set-TextValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  +0.23%  "
Set-TextValue $ws.Range("D37") "0.563"
$ws.Range("E37").Value = "  +4.22%  "
Set-TextValue $ws.Range("D38") "7.69"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").Value = "  -0.75%  "
Set-TextValue $ws.Range("D40") "524.93"
$ws.Range("E40").Value = "  +4.66%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +0.37%  "
Set-TextValue $ws.Range("D43") "0.912"
$ws.Range("E43").Value = "  +10.88%  "
Set-TextValue $ws.Range("D44") "24.12"
$ws.Range("E44").Value = "  -0.94%  "
Set-TextValue $ws.Range("D45") "1.70"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("E49").Value = "  +9.71%  "
Set-TextValue $ws.Range("D50") "53.49"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("E51").Value = "  +2.78%  "
